# VarroaPop Exposed Parameter List - "updated weather handling" edit
#
# Rows 67-75 hold the AI* (Active Ingredient) slope/LD50 parameters that are
# now being generalised to also cover neonicotinoids. The shared explanatory
# note ("Min/Max from Kris' email 8/29/2014") referenced by F67:F75 is
# replaced with an updated note, and the valid min/max bounds on the four
# *Slope / *LD50 rows are widened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updatedNote = "Min/Max from Kris' email 8/29/2014; modified by Jeff to include neonicotinoids"

# AIAdultSlope (row 67): -9..-2 -> 1..9
$ws.Range("D67").Value = 1
$ws.Range("E67").Value = 9
$ws.Range("F67").Value = $updatedNote

# AIAdultLD50 (row 68): 0.01 -> 0 (max of 100 unchanged)
$ws.Range("D68").Value = 0
$ws.Range("F68").Value = $updatedNote

# AIAdultSlopeContact (row 69): -9..-2 -> 1..9
$ws.Range("D69").Value = 1
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = $updatedNote

# AIAdultLD50Contact (row 70): 0.01 -> 0 (max of 100 unchanged)
$ws.Range("D70").Value = 0
$ws.Range("F70").Value = $updatedNote

# AILarvaSlope (row 71): -9..-2 -> 1..9
$ws.Range("D71").Value = 1
$ws.Range("E71").Value = 9
$ws.Range("F71").Value = $updatedNote

# AILarvaLD50 (row 72): 0.01 -> 0 (max of 100 unchanged)
$ws.Range("D72").Value = 0
$ws.Range("F72").Value = $updatedNote

# AIKOW (row 73): min/max unchanged, only the note text is refreshed
$ws.Range("F73").Value = $updatedNote

# AIKOC (row 74): min/max unchanged, only the note text is refreshed
$ws.Range("F74").Value = $updatedNote

# AIHalfLife (row 75): the note cell is removed outright (no footnote here
# anymore) - Clear() drops the <c> element entirely rather than leaving an
# empty, styled cell behind.
$ws.Range("F75").Clear()

# Once no cell references the old shared string "Min/Max from Kris' email
# 8/29/2014" any more, it is pruned automatically from the shared-string
# table on save, and every other cell that referenced a later shared string
# index is renumbered automatically - no further action required there.

# Reflect the author's last on-screen selection/scroll position.
$ws.Activate() | Out-Null
$ws.Range("E72").Select() | Out-Null
